# Update the date line and the multiplication expressions in the table.
$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-08 Sunday", "2024-12-09 Monday"),
    @("988×5=", "675×3="),
    @("824×2=", "219×9="),
    @("623×9=", "189×7="),
    @("567×6=", "427×4="),
    @("885×8=", "625×7="),
    @("839×3=", "313×9="),
    @("696×4=", "117×4="),
    @("366×3=", "439×3="),
    @("671×8=", "596×3="),
    @("766×8=", "645×5="),
    @("418×2=", "946×9="),
    @("761×7=", "875×3="),
    @("263×7=", "657×6="),
    @("619×7=", "727×8="),
    @("526×7=", "302×7="),
    @("850×7=", "390×2="),
    @("247×9=", "372×2="),
    @("450×9=", "627×4="),
    @("869×9=", "152×7="),
    @("873×4=", "884×3="),
    @("908×7=", "935×8="),
    @("194×7=", "395×5="),
    @("194×2=", "581×9="),
    @("122×4=", "242×3="),
    @("743×2=", "815×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
